# CCAA T48 add data for f_e2 -- adds "_f_s" sheets alongside each "_f_e" chart sheet,
# renames cht##_f_e1 -> cht##_f_e, and populates cht11_f_s / cht12_f_e with real data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# cht11: rename cht11_f_e1 -> cht11_f_e, insert cht11_f_s (populated) after it
# ---------------------------------------------------------------------------
$s6 = $wb.Worksheets.Item("cht11_f_e1")

$cht11_f_s = $wb.Worksheets.Add($null, $s6)
$cht11_f_s.Name = "cht11_f_s"
$cht11_f_s.Cells.Item(1, 1).Value = "x"

$s6.Name = "cht11_f_e"
$s6.Cells.Item(1, 2).Value = "f_e"

$cht11_f_s.Cells.Item(1, 2).Value = "f_s"

$colA = @(1,1.0462753950338599,1.0880361173814801,1.1331828442437899,1.1828442437923199,1.2505643340857699,1.3002257336343099,1.3498871331828399,1.4018058690744899,1.4492099322798999,1.5011286681715501,1.5507900677200901,1.6004514672686201,1.6501128668171501,1.6997742663656801,1.7494356659142201,1.80135440180586,1.8510158013544,1.90067720090293,1.95033860045146,2,2.0496613995485302,2.09932279909706,2.1489841986455902,2.1986455981941302,2.2505643340857699,2.2979683972911902,2.3498871331828401,2.3995485327313699,2.4514672686230199,2.5011286681715501,2.5507900677200901,2.6004514672686199,2.6523702031602698,2.6997742663656799,2.75395033860045,2.80135440180586,2.8555304740406302,2.9029345372460398,2.9525959367945802,2.9954853273137698,3)
$colB = @(0.90023014452141803,0.91056433408577797,0.92020990633351496,0.93031447030435099,0.94041799641939705,0.95419553203082397,0.96429905814587002,0.97417269920344496,0.98312628110323996,0.99185101580135404,1.0010344827586199,1.00952881347136,1.0175633740691701,1.0255979346669699,1.03363249526478,1.04097740069017,1.0492413274175501,1.056126462728,1.06278171298098,1.06943696323395,1.0763220985444,1.0827474637399099,1.0887130588204701,1.09467865390104,1.1001844788666599,1.1056897849043801,1.1111961287979,1.1162416647206801,1.12128771957136,1.1261033704366701,1.1309195402298799,1.13504605485067,1.1391725694714701,1.1432985651643699,1.1471957136555799,1.15063153524817,1.15429879868192,1.15773462027451,1.1602524584209,1.1629996626968599,1.1657484237565101,1.1659772709581999)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $cht11_f_s.Cells.Item($r, 1).Value = $colA[$i]
    $cht11_f_s.Cells.Item($r, 2).Value = $colB[$i]
}
$cht11_f_s.Range("A2:B" + ($colA.Length + 1)).NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# cht12: rename cht12_f_e1 -> cht12_f_e and populate it with real data,
# insert cht12_f_s (dummy placeholder) after it
# ---------------------------------------------------------------------------
$s7 = $wb.Worksheets.Item("cht12_f_e1")
$s7.Name = "cht12_f_e"
$s7.Cells.Item(1, 2).Value = "f_e"

$colA2 = @(4.0404040404040398,5.0505050505050502,6.3973063973063899,7.7441077441077404,9.2592592592592595,11.4478114478114,13.973063973063899,16.498316498316498,19.360269360269299,22.5589225589225,25.252525252525199,28.1144781144781,31.481481481481399,35.858585858585798,39.5622895622895,42.929292929292899,45.959595959595902,49.494949494949402,53.872053872053797,58.417508417508401,63.636363636363598,68.855218855218794,74.579124579124496,79.292929292929301,83.501683501683502,86.700336700336706,90.2356902356902,93.771043771043693,97.138047138047099,100.84175084175,104.208754208754,108.080808080808,111.952861952861,116.666666666666,121.717171717171,126.59932659932601,130.47138047138,135.18518518518499,139.89898989898899,144.78114478114401,150)
$colB2 = @(0.801681957186544,0.85382262996941904,0.90091743119265999,0.93960244648318003,0.98165137614678899,1.0321100917431101,1.08425076452599,1.13302752293577,1.17675840978593,1.2204892966360801,1.25581039755351,1.2894495412844,1.3281345565749201,1.3718654434250701,1.4071865443425,1.4357798165137601,1.4593272171253799,1.4862385321100899,1.5114678899082501,1.5333333333333301,1.5501529051987699,1.56360856269113,1.5703363914373001,1.5770642201834799,1.5854740061162,1.5972477064220101,1.61238532110091,1.6308868501528999,1.6561162079510701,1.6880733944954101,1.72003058103975,1.75198776758409,1.7755351681957099,1.80076452599388,1.8243119266054999,1.8411314984709399,1.85290519877675,1.8613149847094801,1.8714067278287401,1.8764525993883701,1.88149847094801)
for ($i = 0; $i -lt $colA2.Length; $i++) {
    $r = $i + 2
    $s7.Cells.Item($r, 1).Value = $colA2[$i]
    $s7.Cells.Item($r, 2).Value = $colB2[$i]
}
$s7.Range("A2:B" + ($colA2.Length + 1)).NumberFormat = "0.000"

$cht12_f_s = $wb.Worksheets.Add($null, $s7)
$cht12_f_s.Name = "cht12_f_s"
$cht12_f_s.Cells.Item(1, 1).Value = "x"
$cht12_f_s.Cells.Item(1, 2).Value = "f_s"
$cht12_f_s.Cells.Item(2, 1).Value = 1
$cht12_f_s.Cells.Item(2, 2).Value = 0.1
$cht12_f_s.Cells.Item(3, 1).Value = 100
$cht12_f_s.Cells.Item(3, 2).Value = 0.999
$cht12_f_s.Range("A2:A3").NumberFormat = "0.000"
$cht12_f_s.Cells.Item(2, 2).NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# cht13: rename cht13_f_e1 -> cht13_f_e (placeholder data untouched),
# insert cht13_f_s (dummy placeholder) after it
# ---------------------------------------------------------------------------
$s8 = $wb.Worksheets.Item("cht13_f_e1")
$s8.Name = "cht13_f_e"
$s8.Cells.Item(1, 2).Value = "f_e"

$cht13_f_s = $wb.Worksheets.Add($null, $s8)
$cht13_f_s.Name = "cht13_f_s"
$cht13_f_s.Cells.Item(1, 1).Value = "x"
$cht13_f_s.Cells.Item(1, 2).Value = "f_s"
$cht13_f_s.Cells.Item(2, 1).Value = 1
$cht13_f_s.Cells.Item(2, 2).Value = 0.1
$cht13_f_s.Cells.Item(3, 1).Value = 100
$cht13_f_s.Cells.Item(3, 2).Value = 0.999
$cht13_f_s.Range("A2:A3").NumberFormat = "0.000"
$cht13_f_s.Cells.Item(2, 2).NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# cht14: rename cht14_f_e1 -> cht14_f_e (placeholder data untouched),
# insert cht14_f_s (dummy placeholder) after it
# ---------------------------------------------------------------------------
$s9 = $wb.Worksheets.Item("cht14_f_e1")
$s9.Name = "cht14_f_e"
$s9.Cells.Item(1, 2).Value = "f_e"

$cht14_f_s = $wb.Worksheets.Add($null, $s9)
$cht14_f_s.Name = "cht14_f_s"
$cht14_f_s.Cells.Item(1, 1).Value = "x"
$cht14_f_s.Cells.Item(1, 2).Value = "f_s"
$cht14_f_s.Cells.Item(2, 1).Value = 1
$cht14_f_s.Cells.Item(2, 2).Value = 0.1
$cht14_f_s.Cells.Item(3, 1).Value = 100
$cht14_f_s.Cells.Item(3, 2).Value = 0.999
$cht14_f_s.Range("A2:A3").NumberFormat = "0.000"
$cht14_f_s.Cells.Item(2, 2).NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# Restore per-sheet selections / view state (best effort)
# ---------------------------------------------------------------------------
$cht14_f_s.Range("B1").Select()
$s9.Range("B1").Select()
$cht13_f_s.Range("B1").Select()
$s8.Range("B1").Select()
$cht12_f_s.Range("B1").Select()
$s7.Range("D3").Select()
$cht11_f_s.Range("B1").Select()

# cht11_f_e ends as the active tab with E28 selected (matches target workbook view)
$s6.Range("E28").Select()

$wb.Save()
